$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "Card-types"

# --- Header cell ---
$ws.Range("A1").Value = "Descriptions"

# --- Card type values (column A) ---
$ws.Range("A2").Value = "Credit Card"
$ws.Range("A3").Value = "Debit Card"
$ws.Range("A4").Value = "Gift Card"
$ws.Range("A5").Value = "Prepaid Card"
$ws.Range("A6").Value = "Store Card"
$ws.Range("A7").Value = "Master Card"
$ws.Range("A8").Value = "Visa Card"

# --- Header alignment: left + vertically centered, keeps wrap ---
$ws.Range("A1").HorizontalAlignment = -4131

# --- Body rows: left + vertically centered ---
$ws.Range("A2:B2").HorizontalAlignment = -4131
$ws.Range("A4:B8").HorizontalAlignment = -4131

# --- Row 3 keeps its own distinct (pre-existing) base style: left only ---
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# --- Trailing blank rows 9,10,12: percentage format + left/center align ---
$ws.Range("A9:B10").HorizontalAlignment = -4131
$ws.Range("A9:B10").NumberFormat = "0%"
$ws.Range("A12:B12").HorizontalAlignment = -4131
$ws.Range("A12:B12").NumberFormat = "0%"

# --- Row 11: plain left/center, same style as the body rows ---
$ws.Range("A11:B11").HorizontalAlignment = -4131

# --- Selection ---
$ws.Range("C7").Select()
